# Add team record (Wins/Losses/Ties) columns to the CIN_1996 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): copy the formatting of the existing header cell (A1,
# bold/centered/bordered) onto the three new header cells, then set text.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Data rows 2-50: every player row gets the team's 1996 record (81-81-0).
for ($r = 2; $r -le 50; $r++) {
    $ws.Cells.Item($r, 29).Value = 81  # AC
    $ws.Cells.Item($r, 30).Value = 81  # AD
    $ws.Cells.Item($r, 31).Value = 0   # AE
}
